$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 09:16"

# Update country stats (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3833597
$ws.Range("C4").Value = 326
$ws.Range("D4").Value = 1775271
$ws.Range("E4").Value = 1915448
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 142878

# Row 6: India
$ws.Range("B6").Value = 1078757
$ws.Range("C6").Value = 893
$ws.Range("D6").Value = 677856
$ws.Range("E6").Value = 374063
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 26838

# Row 53: Armenia
$ws.Range("B53").Value = 34877
$ws.Range("C53").Value = 415
$ws.Range("D53").Value = 23294
$ws.Range("E53").Value = 10942
$ws.Range("G53").Value = 10
$ws.Range("H53").Value = 641

# Row 67: Uzbekistan
$ws.Range("B67").Value = 16429
$ws.Range("C67").Value = 243
$ws.Range("E67").Value = 7219

# Row 75: El Salvador
$ws.Range("D75").Value = 6632
$ws.Range("E75").Value = 4552

# Row 99: Hungria
$ws.Range("B99").Value = 4333
$ws.Range("C99").Value = 18
$ws.Range("D99").Value = 3223
$ws.Range("E99").Value = 514

# Row 139: Letonia
$ws.Range("B139").Value = 1192
$ws.Range("C139").Value = 3
$ws.Range("E139").Value = 139

# Row 146: Georgia
$ws.Range("B146").Value = 1028
$ws.Range("C146").Value = 10
$ws.Range("D146").Value = 899
$ws.Range("E146").Value = 114

# Row 161: Taiwan
$ws.Range("B161").Value = 455
$ws.Range("C161").Value = 1
$ws.Range("E161").Value = 8
